$d = $word.ActiveDocument

# Remove the two trailing "footer" paragraphs that are no longer wanted:
#   "Ver no Jupiter Salvar em pdf Salvar em docx"
#   "© 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution"
# Walk the paragraphs back-to-front (so indices of not-yet-processed
# paragraphs stay valid while we delete) and remove each paragraph (including
# its paragraph mark) whose text matches one of the two footer lines.
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if (($t -like "*Ver no Jupiter*") -or ($t -like "*Powered by Jekyll*")) {
        $p.Range.Delete()
    }
}
